$d = $word.ActiveDocument

# 1. Surat-number line: literal roman-numeral month "II" and literal year
#    "2023" become template placeholders.
$d.Content.Find.Execute("II", $true, $false, $false, $false, $false, $true, 1, $false, "`${bulan_romawi}", 2)
$d.Content.Find.Execute("2023", $true, $false, $false, $false, $false, $true, 1, $false, "`${tahun_surat}", 2)

# 2. Signature block: hard-coded name becomes a template placeholder.
$d.Content.Find.Execute("DWI SAMAYO SATIADY, S.I.K.", $true, $false, $false, $false, $false, $true, 1, $false, "`${kadena}", 2)

# 3. Rank / NRP line: hard-coded rank and NRP become template placeholders.
$d.Content.Find.Execute("KOMISARIS BESAR POLISI NRP 78050947", $true, $false, $false, $false, $false, $true, 1, $false, "`${pangkat_kadena} NRP `${nrp_kadena}", 2)
